$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending ECs -> Target ECs (was Target FAPs)
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.119963
$ws.Range("H2").Value = 3.359889
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.08886233333333333
$ws.Range("N2").Value = 0.266587
$ws.Range("O2").Value = 0.03377239161111666
$ws.Range("P2").Value = 0.03377239161111665
$ws.Range("Q2").Value = 0.09952252542700001
$ws.Range("R2").Value = 0.8957027288430001
$ws.Range("S2").Value = 0.03377239161111666
$ws.Range("T2").Value = 0.03377239161111665

# Row 3: now Target FAPs (was Target sCs)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vip"
$ws.Range("C3").Value = "Adcyap1r1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.119963
$ws.Range("H3").Value = 3.359889
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.124720666666667
$ws.Range("N3").Value = 6.374162
$ws.Range("O3").Value = 0.8075063497346028
$ws.Range("P3").Value = 0.8075063497346027
$ws.Range("Q3").Value = 2.379608532002
$ws.Range("R3").Value = 21.416476788018
$ws.Range("S3").Value = 0.8075063497346028
$ws.Range("T3").Value = 0.8075063497346027

# Row 4: new row, Target sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vip"
$ws.Range("C4").Value = "Adcyap1r1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.119963
$ws.Range("H4").Value = 3.359889
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4176293333333334
$ws.Range("N4").Value = 1.252888
$ws.Range("O4").Value = 0.1587212586542807
$ws.Range("P4").Value = 0.1587212586542807
$ws.Range("Q4").Value = 0.467729401048
$ws.Range("R4").Value = 4.209564609432
$ws.Range("S4").Value = 0.1587212586542807
$ws.Range("T4").Value = 0.1587212586542807
